$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 1.48
$ws.Range("AE2").Value = 8
$ws.Range("G3").Value = 2.8
$ws.Range("I3").Value = 2.75
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 1.7
$ws.Range("U3").Value = 12
$ws.Range("V3").Value = 11
$ws.Range("X3").Value = 26
$ws.Range("AB3").Value = 17
$ws.Range("AE3").Value = 7
$ws.Range("AF3").Value = 12
$ws.Range("AH3").Value = 29
$ws.Range("G7").Value = 2
$ws.Range("I7").Value = 4.33
$ws.Range("J7").Value = 1.13
$ws.Range("K7").Value = 6
$ws.Range("V7").Value = 9.5
$ws.Range("AH7").Value = 51
$ws.Range("P8").Value = 1.58
$ws.Range("P9").Value = 1.58
$ws.Range("P10").Value = 1.63
$ws.Range("H12").Value = 3.4
$ws.Range("L12").Value = 1.33
$ws.Range("M12").Value = 3.25
$ws.Range("N12").Value = 2.08
$ws.Range("O12").Value = 1.73
$ws.Range("P12").Value = 1.44
$ws.Range("Q12").Value = 2.63
$ws.Range("T12").Value = 6.5
$ws.Range("U12").Value = 8
$ws.Range("Y12").Value = 29
$ws.Range("Z12").Value = 9
$ws.Range("AA12").Value = 6.5
$ws.Range("AB12").Value = 17
$ws.Range("AD12").Value = 301
$ws.Range("AF12").Value = 21
$ws.Range("AI12").Value = 41
$ws.Range("G13").Value = 3.8
$ws.Range("I13").Value = 2
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.5
$ws.Range("T13").Value = 9.5
$ws.Range("AF13").Value = 9
$ws.Range("AJ13").Value = 34
$ws.Range("J14").Value = 1.08
$ws.Range("K14").Value = 8
$ws.Range("N14").Value = 2.15
$ws.Range("O14").Value = 1.67
$ws.Range("L15").Value = 1.29
$ws.Range("M15").Value = 3.5
$ws.Range("N15").Value = 1.93
$ws.Range("O15").Value = 1.93
$ws.Range("J22").Value = 1.11
$ws.Range("K22").Value = 5.5
$ws.Range("L22").Value = 1.53
$ws.Range("N22").Value = 2.55
$ws.Range("Q22").Value = 2.22
$ws.Range("Z22").Value = 5.5
$ws.Range("AH22").Value = 80
$ws.Range("G23").Value = 1.4
$ws.Range("H23").Value = 3.95
$ws.Range("I23").Value = 7.7
$ws.Range("M23").Value = 2.87
$ws.Range("N23").Value = 2.05
$ws.Range("O23").Value = 1.7
$ws.Range("T23").Value = 5.2
$ws.Range("X23").Value = 13.5
$ws.Range("AC23").Value = 200
$ws.Range("AE23").Value = 15.5
$ws.Range("AF23").Value = 50
$ws.Range("AG23").Value = 27
$ws.Range("G24").Value = 2.55
$ws.Range("I24").Value = 2.57
$ws.Range("L24").Value = 1.39
$ws.Range("M24").Value = 2.77
$ws.Range("P24").Value = 1.5
$ws.Range("Q24").Value = 2.42
$ws.Range("R24").Value = 1.88
$ws.Range("S24").Value = 1.82
$ws.Range("T24").Value = 7.5
$ws.Range("W24").Value = 28
$ws.Range("X24").Value = 23
$ws.Range("AE24").Value = 7.4
$ws.Range("AF24").Value = 12
$ws.Range("AI24").Value = 24
$ws.Range("AJ24").Value = 37
$ws.Range("G26").Value = 2.35
$ws.Range("I26").Value = 3.25
$ws.Range("O26").Value = 1.5
$ws.Range("T26").Value = 6.5
$ws.Range("U26").Value = 10
$ws.Range("W26").Value = 23
$ws.Range("X26").Value = 21
$ws.Range("AE26").Value = 8
$ws.Range("AF26").Value = 15
$ws.Range("AH26").Value = 34
$ws.Range("J27").Value = 1.06
$ws.Range("K27").Value = 10
$ws.Range("N27").Value = 2.05
$ws.Range("O27").Value = 1.72
$ws.Range("N28").Value = 1.62
$ws.Range("O28").Value = 2.2
$ws.Range("U28").Value = 7
$ws.Range("N30").Value = 2.07
$ws.Range("O30").Value = 1.69
$ws.Range("J31").Value = 1.07
$ws.Range("K31").Value = 7
$ws.Range("L31").Value = 1.31
$ws.Range("M31").Value = 3.15
$ws.Range("N31").Value = 1.93
$ws.Range("O31").Value = 1.78
$ws.Range("P31").Value = 1.44
$ws.Range("Q31").Value = 2.57
$ws.Range("R31").Value = 1.75
$ws.Range("S31").Value = 1.98
$ws.Range("T31").Value = 7.2
$ws.Range("U31").Value = 9.5
$ws.Range("W31").Value = 18
$ws.Range("X31").Value = 15.5
$ws.Range("Z31").Value = 7
$ws.Range("AA31").Value = 6.2
$ws.Range("AB31").Value = 14
$ws.Range("AC31").Value = 65
$ws.Range("AD31").Value = 500
$ws.Range("AE31").Value = 10.5
$ws.Range("AF31").Value = 21
$ws.Range("AG31").Value = 12.5
$ws.Range("AI31").Value = 35
$ws.Range("AJ31").Value = 40
$ws.Range("AE32").Value = 5.5
$ws.Range("AG32").Value = 8.5
$ws.Range("Q33").Value = 2.75
$ws.Range("R33").Value = 1.65
$ws.Range("S33").Value = 2.12
$ws.Range("T33").Value = 8.5
$ws.Range("X33").Value = 14
$ws.Range("AE33").Value = 11.5
$ws.Range("AF33").Value = 21
$ws.Range("AH33").Value = 55
$ws.Range("AJ33").Value = 37
$ws.Range("H34").Value = 3.95
$ws.Range("I34").Value = 4.9
$ws.Range("L34").Value = 1.24
$ws.Range("M34").Value = 3.65
$ws.Range("N34").Value = 1.72
$ws.Range("Q34").Value = 2.92
$ws.Range("T34").Value = 7.4
$ws.Range("V34").Value = 8
$ws.Range("AA34").Value = 7.9
$ws.Range("AF34").Value = 30
$ws.Range("AG34").Value = 16.5
$ws.Range("AI34").Value = 50
$ws.Range("G35").Value = 1.4
$ws.Range("I35").Value = 7.5
$ws.Range("J35").Value = 1.08
$ws.Range("K35").Value = 8
$ws.Range("N35").Value = 2.3
$ws.Range("O35").Value = 1.6
$ws.Range("R35").Value = 2.75
$ws.Range("S35").Value = 1.37
$ws.Range("T35").Value = 4.75
$ws.Range("AB35").Value = 34
$ws.Range("AC35").Value = 151
$ws.Range("AF35").Value = 41
$ws.Range("AJ35").Value = 81
$ws.Range("R36").Value = 1.87
$ws.Range("S36").Value = 1.87
$ws.Range("R37").Value = 1.77
$ws.Range("S37").Value = 1.87
$ws.Range("G38").Value = 3.6
$ws.Range("I38").Value = 1.95
$ws.Range("J38").Value = 1.03
$ws.Range("K38").Value = 15
$ws.Range("R38").Value = 1.54
$ws.Range("U38").Value = 21
$ws.Range("V38").Value = 13
$ws.Range("AH38").Value = 17
$ws.Range("G39").Value = 2.15
$ws.Range("H39").Value = 2.8
$ws.Range("I39").Value = 3.8
$ws.Range("M39").Value = 2
$ws.Range("N39").Value = 2.8
$ws.Range("O39").Value = 1.33
$ws.Range("P39").Value = 1.62
$ws.Range("Q39").Value = 2.02
$ws.Range("R39").Value = 2.32
$ws.Range("S39").Value = 1.47
$ws.Range("T39").Value = 4.9
$ws.Range("U39").Value = 8.25
$ws.Range("V39").Value = 10.25
$ws.Range("X39").Value = 25
$ws.Range("Z39").Value = 4.5
$ws.Range("AA39").Value = 5.9
$ws.Range("AB39").Value = 23
$ws.Range("AE39").Value = 7.1
$ws.Range("AF39").Value = 17.5
$ws.Range("AH39").Value = 65
$ws.Range("AI39").Value = 55
$ws.Range("G40").Value = 2.2
$ws.Range("H40").Value = 2.92
$ws.Range("I40").Value = 3.4
$ws.Range("K40").Value = 6
$ws.Range("L40").Value = 1.45
$ws.Range("M40").Value = 2.37
$ws.Range("N40").Value = 2.27
$ws.Range("O40").Value = 1.5
$ws.Range("P40").Value = 1.5
$ws.Range("Q40").Value = 2.27
$ws.Range("R40").Value = 1.98
$ws.Range("S40").Value = 1.65
$ws.Range("T40").Value = 6
$ws.Range("U40").Value = 9.5
$ws.Range("V40").Value = 9.25
$ws.Range("W40").Value = 21
$ws.Range("X40").Value = 21
$ws.Range("Y40").Value = 40
$ws.Range("Z40").Value = 6.7
$ws.Range("AB40").Value = 17
$ws.Range("AC40").Value = 100
$ws.Range("AE40").Value = 7.9
$ws.Range("AF40").Value = 16.5
$ws.Range("AG40").Value = 12.5
$ws.Range("AH40").Value = 50
$ws.Range("AJ40").Value = 50
$ws.Range("G41").Value = 3.3
$ws.Range("H41").Value = 3.25
$ws.Range("I41").Value = 2.1
$ws.Range("J41").Value = 1.06
$ws.Range("K41").Value = 10
$ws.Range("L41").Value = 1.3
$ws.Range("M41").Value = 3.4
$ws.Range("N41").Value = 2.03
$ws.Range("O41").Value = 1.83
$ws.Range("P41").Value = 1.44
$ws.Range("Q41").Value = 2.63
$ws.Range("R41").Value = 1.8
$ws.Range("S41").Value = 1.91
$ws.Range("T41").Value = 10
$ws.Range("U41").Value = 17
$ws.Range("W41").Value = 34
$ws.Range("Z41").Value = 9.5
$ws.Range("AD41").Value = 251
$ws.Range("AE41").Value = 7.5
$ws.Range("AJ41").Value = 29
